# Update of all values to match PDF edition 10 (commit 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# segment | y2019 | share_of_total_2019 | y2021 | share_of_total_2021 | y2022 | share_of_total_2022
$ws.Range("A1").Value = "segment"
$ws.Range("B1").Value = "y2019"
$ws.Range("C1").Value = "share_of_total_2019"
$ws.Range("D1").Value = "y2021"
$ws.Range("E1").Value = "share_of_total_2021"
$ws.Range("F1").Value = "y2022"
$ws.Range("G1").Value = "share_of_total_2022"

# ---- Row 2: Mainline ----
$ws.Range("A2").Value = "Mainline"
$ws.Range("B2").Value = 3991685
$ws.Range("C2").Value = 0.36
$ws.Range("D2").Value = 1816909
$ws.Range("E2").Value = 0.29
$ws.Range("F2").Value = 2981880
$ws.Range("G2").Value = 0.32

# ---- Row 3: Low-cost ----
$ws.Range("A3").Value = "Low-cost"
$ws.Range("B3").Value = 3493913
$ws.Range("C3").Value = 0.32
$ws.Range("D3").Value = 1610239
$ws.Range("E3").Value = 0.26
$ws.Range("F3").Value = 2984376
$ws.Range("G3").Value = 0.32

# ---- Row 4: Regional ----
$ws.Range("A4").Value = "Regional"
$ws.Range("B4").Value = 1643854
$ws.Range("C4").Value = 0.15
$ws.Range("D4").Value = 861587
$ws.Range("E4").Value = 0.14
$ws.Range("F4").Value = 1219685
$ws.Range("G4").Value = 0.13

# ---- Row 5: Business Aviation ----
$ws.Range("A5").Value = "Business Aviation"
$ws.Range("B5").Value = 683473
$ws.Range("C5").Value = 0.06
$ws.Range("D5").Value = 709398
$ws.Range("E5").Value = 0.11
$ws.Range("F5").Value = 791909
$ws.Range("G5").Value = 0.09

# ---- Row 6: All-Cargo ----
$ws.Range("A6").Value = "All-Cargo"
$ws.Range("B6").Value = 368362
$ws.Range("C6").Value = 0.03
$ws.Range("D6").Value = 419824
$ws.Range("E6").Value = 0.07
$ws.Range("F6").Value = 389611
$ws.Range("G6").Value = 0.04

# ---- Row 7: Other ----
$ws.Range("A7").Value = "Other"
$ws.Range("B7").Value = 372796
$ws.Range("C7").Value = 0.03
$ws.Range("D7").Value = 363712
$ws.Range("E7").Value = 0.06
$ws.Range("F7").Value = 389396
$ws.Range("G7").Value = 0.04

# ---- Row 8: Charter ----
$ws.Range("A8").Value = "Charter"
$ws.Range("B8").Value = 382218
$ws.Range("C8").Value = 0.04
$ws.Range("D8").Value = 303384
$ws.Range("E8").Value = 0.05
$ws.Range("F8").Value = 324824
$ws.Range("G8").Value = 0.04

# ---- Row 9: Military ----
$ws.Range("A9").Value = "Military"
$ws.Range("B9").Value = 149001
$ws.Range("C9").Value = 0.01
$ws.Range("D9").Value = 145699
$ws.Range("E9").Value = 0.02
$ws.Range("F9").Value = 156012
$ws.Range("G9").Value = 0.02

# ---- Row 10: Total ----
$ws.Range("A10").Value = "Total"
$ws.Range("B10").Value = 11085302
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 6230752
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 9237693
$ws.Range("G10").Value = 1

# ---- New blank cells below the table carrying the Percent style ----
$cells = @("D12", "F12", "B13", "B14", "B15", "D16")
foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Percent"
}

# ---- Selection matches the saved view state ----
$ws.Range("A10").Select()
